$d = $word.ActiveDocument
$br = [char]11

# ---------------------------------------------------------------
# 1. Global font-name fix: "TimesNewToman" -> "Times New Roman"
#    (applies to every run already in the document)
# ---------------------------------------------------------------
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"

# ---------------------------------------------------------------
# 2. Title / author / e-mail replacements
# ---------------------------------------------------------------
$d.Content.Find.Execute("Visions of the Quantum Realm", $false, $false, $false, $false, $false, $true, 1, $false, "The Enchanting Symphony of Colors", 2) | Out-Null
$d.Content.Find.Execute("Charlene Evans", $false, $false, $false, $false, $false, $true, 1, $false, "Roselyn Carter", 2) | Out-Null
$d.Content.Find.Execute("cevans@wordmail", $false, $false, $false, $false, $false, $true, 1, $false, "roselyncarter@xyzschool", 2) | Out-Null
$d.Content.Find.Execute("net", $false, $true, $false, $false, $false, $true, 1, $false, "edu", 2) | Out-Null

# ---------------------------------------------------------------
# 3. Rewrite the big body paragraph (paragraph 5) completely
# ---------------------------------------------------------------
function Insert-Token {
    param($para, $text)
    $r = $para.Range
    $r.End = $r.End - 1
    $startPos = $r.End
    $r.InsertAfter($text)
    $endPos = $startPos + $text.Length
    $newRange = $d.Range($startPos, $endPos)
    $newRange.Font.Name = "Times New Roman"
    $newRange.Font.Size = 12
    $newRange.Font.Color = 0
    return $newRange
}

$bodyPara = $d.Paragraphs.Item(5)
$bodyRange = $bodyPara.Range
$bodyRange.End = $bodyRange.End - 1
$bodyRange.Text = ""

$bodyTokens = @(
    "- Dive into the vibrant realm of colors, where hues dance together in a captivating symphony",
    ".",
    " From the blazing scarlet sunsets that ignite the skies to the calming cerulean depths of the ocean, colors enchant our world with their unspoken stories",
    ".",
    "${br}",
    "${br}- Colors, like musical notes, possess the power to evoke emotions and shape perceptions",
    ".",
    " From the vibrant reds that spark passion to the soothing greens that promote tranquility, colors have a profound impact on our psychological and physiological well-being",
    ".",
    "${br}",
    "${br}- Beyond their aesthetic allure, colors play a crucial role in various scientific fields",
    ".",
    " Whether it's the study of light and its interactions or the analysis of chemical compounds, colors serve as essential tools for unraveling the mysteries of the universe",
    ".",
    "${br}",
    "${br}Body:",
    "${br}",
    "${br}Paragraph 1:",
    "${br}",
    "${br}- In the realm of art, colors become the language of expression",
    ".",
    " Artists use colors to convey emotions, tell stories, and create visual masterpieces",
    ".",
    " From the bold strokes of abstract paintings to the intricate details of realistic landscapes, colors allow artists to share their unique perspectives and connect with viewers on a profound level",
    ".",
    "${br}",
    "${br}Paragraph 2:",
    "${br}",
    "${br}- In the realm of science, colors hold significant importance",
    ".",
    " From the rainbow's spectrum, scientists gain insights into the properties of light and its interactions with matter",
    ".",
    " Colors also serve as indicators in chemical reactions, revealing the composition and structure of substances",
    ".",
    " Moreover, colors play a vital role in biotechnology and medical research, assisting in the development of diagnostic techniques and treatments",
    ".",
    "${br}",
    "${br}Paragraph 3:",
    "${br}",
    "${br}- In the realm of history and culture, colors carry immense significance",
    ".",
    " Different cultures associate specific colors with emotions, values, and beliefs",
    ".",
    " From the vibrant ",
    "colors of traditional festivals to the symbolic use of colors in flags and emblems, colors serve as cultural markers that reflect the heritage and identity of nations",
    ".",
    " Colors also play a crucial role in politics, influencing public perception and electoral outcomes",
    "."
)

foreach ($tok in $bodyTokens) {
    Insert-Token $bodyPara $tok | Out-Null
}

# ---------------------------------------------------------------
# 4. Rewrite the Summary content paragraph (paragraph 7) sentence by
#    sentence so the existing run/formatting boundaries are kept.
# ---------------------------------------------------------------
$d.Content.Find.Execute("The quantum realm, with its enigmatic phenomena like superposition and entanglement, challenges our classical understanding of reality", $false, $false, $false, $false, $false, $true, 1, $false, "- Colors, with their inherent beauty and symbolism, hold a profound impact on our lives", 2) | Out-Null
$d.Content.Find.Execute(" Its exploration has unveiled a hidden interconnectedness at the heart of the universe, with far-reaching implications across disciplines", $false, $false, $false, $false, $false, $true, 1, $false, " From their role in art and science to their significance in culture and history, colors shape our perceptions, evoke emotions, and contribute to our understanding of the world around us", 2) | Out-Null
$d.Content.Find.Execute(" As we continue to unravel the mysteries of quantum mechanics, we open up new avenues for innovation and discovery, transforming our understanding of the world and our place within it", $false, $false, $false, $false, $false, $true, 1, $false, " By delving into the fascinating world of colors, we unlock a treasure trove of knowledge, beauty, and inspiration", 2) | Out-Null

# ---------------------------------------------------------------
# 5. Append a new empty paragraph at the very end of the document
# ---------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Host "Edit complete."
